# olcsaprojektdoku.xlsx - add two new "Feladat" rows (4 and 5) to Munka1,
# one each for the two existing team members, with a later date and new
# task descriptions, then tidy up column widths / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---- Row 4: Heizer Marcell ------------------------------------------------
$ws.Range("A4").Value = "Heizer Marcell"

$ws.Range("C4").Value = 45611
$ws.Range("C2").Copy()
$ws.Range("C4:D4").PasteSpecial($xlPasteFormats)

$ws.Range("E4").Value = "Kód befejezése, a dokumentációban képekért felelős"
$ws.Range("E2").Copy()
$ws.Range("E4:F4").PasteSpecial($xlPasteFormats)

$null = $ws.Range("C4:D4").Merge()
$null = $ws.Range("E4:F4").Merge()

# ---- Row 5: Benedek Zsombor ------------------------------------------------
$ws.Range("A5").Value = "Benedek Zsombor"

$ws.Range("C5").Value = 45611
$ws.Range("C3").Copy()
$ws.Range("C5:D5").PasteSpecial($xlPasteFormats)

$ws.Range("E5").Value = "Kód befejezése, dokumentáció szövegének elkészítése"
$ws.Range("E3").Copy()
$ws.Range("E5:F5").PasteSpecial($xlPasteFormats)

$null = $ws.Range("C5:D5").Merge()
$null = $ws.Range("E5:F5").Merge()

$excel.CutCopyMode = $false

# ---- column widths ----------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 9.25
$ws.Columns.Item(6).ColumnWidth = 40.26

# ---- selection / active cell -------------------------------------------------
$null = $ws.Range("I5").Select()
